{"js": "// Replace the date line and the 25 \"three-digit \u00f7 one-digit\" problems with\n// their updated values, in document order. Several of the old/new values\n// collide (e.g. \"954\u00f77=\" is both a target in one cell and a source in\n// another, \"438\u00f73=\" likewise), so the replacement is done positionally\n// (by paragraph order) rather than by a global text search/replace, which\n// would risk re-matching an already-updated cell.\nconst replacements = [\n  \"2024-04-25 Thursday\",\n  \"954\u00f77=\",\n  \"946\u00f78=\",\n  \"703\u00f79=\",\n  \"526\u00f75=\",\n  \"517\u00f72=\",\n  \"648\u00f73=\",\n  \"307\u00f77=\",\n  \"512\u00f72=\",\n  \"275\u00f79=\",\n  \"438\u00f73=\",\n  \"529\u00f73=\",\n  \"280\u00f76=\",\n  \"444\u00f76=\",\n  \"584\u00f74=\",\n  \"719\u00f78=\",\n  \"690\u00f74=\",\n  \"636\u00f72=\",\n  \"113\u00f78=\",\n  \"368\u00f78=\",\n  \"631\u00f75=\",\n  \"433\u00f76=\",\n  \"563\u00f76=\",\n  \"430\u00f76=\",\n  \"257\u00f79=\",\n  \"638\u00f79=\",\n];\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\n// Only paragraphs that actually hold text are in play (the table also has\n// several intentionally blank rows/cells that must stay untouched).\nlet repIdx = 0;\nfor (let i = 0; i < paragraphs.items.length && repIdx < replacements.length; i++) {\n  const para = paragraphs.items[i];\n  if (para.text && para.text.length > 0) {\n    para.getRange().insertText(replacements[repIdx], Word.InsertLocation.replace);\n    repIdx++;\n  }\n}\n\nawait context.sync();\n", "ps1": "# Replace the date line and the 25 \"three-digit \u00f7 one-digit\" problems with\n# their updated values. The problems live in a 20-row x 5-column table where\n# every 4th row (1, 5, 9, 13, 17) holds data and the rows in between are\n# intentionally blank, so addressing cells via Table.Cell(row, col) is the\n# safest way to hit exactly the right run -- several of the old/new values\n# collide (e.g. \"954\u00f77=\" is both a target in one cell and a source in\n# another, \"438\u00f73=\" likewise), so a positional update (not a global\n# Find/Replace) is required to avoid re-matching an already-updated cell.\n$d = $word.ActiveDocument\n\n# Update the date line (first paragraph, above the table).\n$dateRange = $d.Paragraphs.Item(1).Range\n[void]$dateRange.MoveEnd(1, -1)\n$dateRange.Text = \"2024-04-25 Thursday\"\n\n$t = $d.Tables.Item(1)\n$dataRows = @(1, 5, 9, 13, 17)\n$values = @(\n    @(\"954\u00f77=\", \"946\u00f78=\", \"703\u00f79=\", \"526\u00f75=\", \"517\u00f72=\"),\n    @(\"648\u00f73=\", \"307\u00f77=\", \"512\u00f72=\", \"275\u00f79=\", \"438\u00f73=\"),\n    @(\"529\u00f73=\", \"280\u00f76=\", \"444\u00f76=\", \"584\u00f74=\", \"719\u00f78=\"),\n    @(\"690\u00f74=\", \"636\u00f72=\", \"113\u00f78=\", \"368\u00f78=\", \"631\u00f75=\"),\n    @(\"433\u00f76=\", \"563\u00f76=\", \"430\u00f76=\", \"257\u00f79=\", \"638\u00f79=\")\n)\n\nfor ($ri = 0; $ri -lt $dataRows.Length; $ri++) {\n    $row = $dataRows[$ri]\n    for ($c = 1; $c -le 5; $c++) {\n        $cell = $t.Cell($row, $c)\n        $r = $cell.Range\n        # Drop the trailing paragraph mark + cell mark so only the cell's\n        # real text is overwritten. wdCharacter == 1.\n        [void]$r.MoveEnd(1, -1)\n        $r.Text = $values[$ri][$c - 1]\n    }\n}\n"}
